$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-31 Wednesday" "2024-02-01 Thursday"

Replace-Text "13×81=" "78×52="
Replace-Text "92×38=" "32×73="
Replace-Text "62×81=" "54×77="
Replace-Text "92×83=" "50×21="
Replace-Text "42×56=" "54×13="
Replace-Text "85×56=" "96×32="
Replace-Text "44×18=" "18×31="
Replace-Text "43×71=" "21×35="
Replace-Text "32×94=" "98×38="
Replace-Text "41×63=" "40×63="
Replace-Text "17×85=" "11×60="
Replace-Text "51×87=" "27×42="
Replace-Text "86×14=" "48×67="
Replace-Text "99×78=" "64×85="
Replace-Text "29×56=" "72×67="
Replace-Text "50×86=" "31×60="
Replace-Text "44×61=" "81×45="
Replace-Text "35×44=" "96×69="
Replace-Text "65×83=" "66×75="
Replace-Text "19×86=" "33×78="
Replace-Text "19×84=" "33×28="
Replace-Text "56×63=" "64×71="
Replace-Text "26×32=" "60×30="
Replace-Text "43×85=" "39×88="
Replace-Text "17×48=" "60×73="
